$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: PartnerName -> Vendor
$ws.Range("A1").Value = "Vendor"

# Update "Ad4Game" -> "Ad 4Game" in column A (rows 2-5)
$ws.Range("A2").Value = "Ad 4Game"
$ws.Range("A3").Value = "Ad 4Game"
$ws.Range("A4").Value = "Ad 4Game"
$ws.Range("A5").Value = "Ad 4Game"

# Update selection to match author's final cursor position
$ws.Range("C5").Select()
